$d = $word.ActiveDocument

# 1) Remove the old _GoBack bookmark (will be re-added later in the "Slide 30" paragraph)
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2) Rewrite the "Slide 28-29" paragraph (most common words / quoted skill lists)
$para24xml = '<w:p><w:pPr><w:spacing w:after="0" w:line="276" w:lineRule="auto"/></w:pPr><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t xml:space="preserve">Slide </w:t></w:r><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>28</w:t></w:r><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>-29</w:t></w:r><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>:</w:t></w:r><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Now that we had our analysis complete, we went back to our</w:t></w:r><w:r><w:t xml:space="preserve"> job postings </w:t></w:r><w:r><w:t>to pull out the most common word</w:t></w:r><w:r><w:t>s</w:t></w:r><w:r><w:t xml:space="preserve">. </w:t></w:r><w:r><w:t>Here you can see</w:t></w:r><w:r><w:t xml:space="preserve"> that Data Scientist, Data Analyst, and Statistician all share common words like “Analysis”, “Statistics”</w:t></w:r><w:r><w:t>, “Communication”, and “Team”</w:t></w:r><w:r><w:t xml:space="preserve">. A Data Scientist adds tools </w:t></w:r><w:r><w:t xml:space="preserve">and skills </w:t></w:r><w:r><w:t xml:space="preserve">such as </w:t></w:r><w:r><w:t>“</w:t></w:r><w:r><w:t>Programming</w:t></w:r><w:r><w:t>”</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>“</w:t></w:r><w:r><w:t>Engineering</w:t></w:r><w:r><w:t>”</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t xml:space="preserve">and </w:t></w:r><w:r><w:t>“</w:t></w:r><w:r><w:t>Machine Learning</w:t></w:r><w:r><w:t>”</w:t></w:r><w:r><w:t xml:space="preserve">. </w:t></w:r></w:p>'
$xml24 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $para24xml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$found24 = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*pull out the most common words*") {
        $d.Paragraphs($i).Range.InsertXML($xml24)
        $found24 = $true
        break
    }
}
if (-not $found24) {
    throw "Could not locate the 'Slide 28-29' paragraph"
}

# 3) Rewrite the "Slide 30" paragraph (Personal Takeaway) and re-add the _GoBack bookmark
$para26xml = '<w:p><w:pPr><w:spacing w:line="276" w:lineRule="auto"/></w:pPr><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>S</w:t></w:r><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t xml:space="preserve">lide </w:t></w:r><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>30</w:t></w:r><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r><w:r><w:t>Personal Takeaway</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
$xml26 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $para26xml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$found26 = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*When considering*") {
        $d.Paragraphs($i).Range.InsertXML($xml26)
        $found26 = $true
        break
    }
}
if (-not $found26) {
    throw "Could not locate the 'Slide 30' paragraph"
}
